# PV Default Costs.xlsx - "Proposed 2016.3.14" sheet update
# Adds a new "Utility / Proposed" (G) data column alongside the existing
# "Utility / Existing" (F) column, fills in the "Proposed" (C/E/G) values
# for Residential/Commercial/Utility that were previously blank, renames a
# couple of row labels (with a yellow highlight to flag the rename), and
# extends the weighted-total formulas in row 18 to the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proposed 2016.3.14")

# ---------------------------------------------------------------------
# 1. Row-label text changes (order matters: it fixes the order new shared
#    strings are appended in, matching the source file).
# ---------------------------------------------------------------------
$ws.Range("A12").Value = 'Engineering ($/Wdc): Changed to be "Developer Overhead ($/Wdc)"'
$ws.Range("A15").Value = 'Land preparation: Changed to be "Transmission Line ($/Wdc)"'
$ws.Range("A5").Value  = "Inverter (`$/Wdc)"

# Highlight the two renamed labels in yellow.
$ws.Range("A12").Interior.Color = 65535
$ws.Range("A15").Interior.Color = 65535

# ---------------------------------------------------------------------
# 2. Residential / Commercial "Proposed" column (C) - previously blank.
# ---------------------------------------------------------------------
$ws.Range("C4").Value = 0.69666666666666699
$ws.Range("C5").Value = 0.28999999999999998
$ws.Range("C6").Formula = "=0.329105205267494+0.203586501897641"
$ws.Range("C7").Formula = "=0.215199265236845+0.115106555743999"
$ws.Range("C8").Formula = "=0.677401030654478+0.37787140321983"
$ws.Range("C9").Value = 0
$ws.Range("C11").Value = 0.11
$ws.Range("C16").Value = 0.51679868381138649
$ws.Range("C17").Value = 0.05
$ws.Range("C18").Formula = "=SUM(C4:C8)*(1+C9)*(1+C16*C17)+SUM(C11:C15)"

# ---------------------------------------------------------------------
# 3. Commercial "Proposed" column (E) - previously blank.
# ---------------------------------------------------------------------
$ws.Range("E4").Value = 0.68
$ws.Range("E5").Value = 0.13
$ws.Range("E6").Formula = "=0.178+0.156+0.029"
$ws.Range("E7").Formula = "=0.55*0.187+0.45*0.187"
$ws.Range("E8").Formula = "=0.06+0.1661/2+0.1661/2+0.493"
$ws.Range("E9").Value = 0
$ws.Range("E11").Value = 0.001
$ws.Range("E16").Value = 0.67240632966186697
$ws.Range("E17").Value = 0.05
$ws.Range("E18").Formula = "=SUM(E4:E8)*(1+E9)*(1+E16*E17)+SUM(E11:E15)"

# ---------------------------------------------------------------------
# 4. New Utility "Proposed" column (G) - new column, mirrors F (Existing).
# ---------------------------------------------------------------------
$ws.Range("G4").Value = 0.65
$ws.Range("G5").Value = 0.11
$ws.Range("G6").Formula = "=0.16+0.16"
$ws.Range("G7").Value = 0.19
$ws.Range("G8").Formula = "=0.103724570767137+0.06"
$ws.Range("G9").Value = 0
$ws.Range("G11").Value = 0.023767500513328502
$ws.Range("G12").Formula = "=0.16"
$ws.Range("G13").Value = 0.03
$ws.Range("G14").Value = 0.03
$ws.Range("G15").Value = 0.02
$ws.Range("G16").Value = 1
$ws.Range("G17").Value = 0.05
$ws.Range("G18").Formula = "=SUM(G4:G8)*(1+G9)*(1+G16*G17)+SUM(G11:G15)"

# Make F13/F14/F15 explicit values too (style refresh below needs them).
$ws.Range("F13").Value = 0.03
$ws.Range("F14").Value = 0.03
$ws.Range("F15").Value = 0.06

# ---------------------------------------------------------------------
# 5. Currency number format ("$"#,##0.00) on the Utility (F/G) value
#    cells, matching the rest of the "Existing" column's look.
# ---------------------------------------------------------------------
$ws.Range("F4:G8").NumberFormat = '"$"#,##0.00'
$ws.Range("F11:G15").NumberFormat = '"$"#,##0.00'

# ---------------------------------------------------------------------
# 6. Row heights (rows already at 14.45pt get an explicit custom height).
# ---------------------------------------------------------------------
$heightRows = @(4,5,6,7,8,9,10,11,16)
foreach ($r in $heightRows) {
    $ws.Rows($r).RowHeight = 14.45
}

# ---------------------------------------------------------------------
# 7. Selection moves to A5 (matches the saved view state in the source).
# ---------------------------------------------------------------------
$ws.Range("A5").Select() | Out-Null
